$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.021026094879971
$ws.Range("D2").Value = 1.036664810600258
$ws.Range("E2").Value = 1.021959839254618
$ws.Range("F2").Value = 1.03403738294863
$ws.Range("I2").Value = 1.035886513026745
$ws.Range("J2").Value = 1.02621988583449
$ws.Range("K2").Value = 1.039457903987421
$ws.Range("L2").Value = 1.024795630086122
$ws.Range("M2").Value = 1.036838008525988
$ws.Range("N2").Value = 1.012816112434217
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.021854791337615
$ws.Range("D3").Value = 1.037219448476655
$ws.Range("E3").Value = 1.022659515899911
$ws.Range("F3").Value = 1.035133808882917
$ws.Range("I3").Value = 1.036089687960619
$ws.Range("J3").Value = 1.026686652994488
$ws.Range("K3").Value = 1.039822866599263
$ws.Range("L3").Value = 1.025302163328352
$ws.Range("M3").Value = 1.037742773187417
$ws.Range("N3").Value = 1.012971208126736
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.022391422593479
$ws.Range("D4").Value = 1.037578184431018
$ws.Range("E4").Value = 1.023113007619875
$ws.Range("F4").Value = 1.035843857112496
$ws.Range("I4").Value = 1.036219803321524
$ws.Range("J4").Value = 1.026988451592246
$ws.Range("K4").Value = 1.040058164602893
$ws.Range("L4").Value = 1.025630005086323
$ws.Range("M4").Value = 1.038328216743912
$ws.Range("N4").Value = 1.013071460799079
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.022617119026455
$ws.Range("D5").Value = 1.037728958888314
$ws.Range("E5").Value = 1.023303834525644
$ws.Range("F5").Value = 1.036142501282212
$ws.Range("I5").Value = 1.036274179272664
$ws.Range("J5").Value = 1.027115271431206
$ws.Range("K5").Value = 1.040156877391249
$ws.Range("L5").Value = 1.025767847997382
$ws.Range("M5").Value = 1.03857433652289
$ws.Range("N5").Value = 1.013113581610843
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.022655020080138
$ws.Range("D6").Value = 1.037754272294609
$ws.Range("E6").Value = 1.023335885653338
$ws.Range("F6").Value = 1.036192653135156
$ws.Range("I6").Value = 1.036283290179931
$ws.Range("J6").Value = 1.027136561710211
$ws.Range("K6").Value = 1.040173439556827
$ws.Range("L6").Value = 1.025790993465263
$ws.Range("M6").Value = 1.038615661053875
$ws.Range("N6").Value = 1.01312065237963
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.022394437982131
$ws.Range("D7").Value = 1.037580199238811
$ws.Range("E7").Value = 1.023115556756996
$ws.Range("F7").Value = 1.035847847061711
$ws.Range("I7").Value = 1.036220531171866
$ws.Range("J7").Value = 1.026990146387581
$ws.Range("K7").Value = 1.040059484420983
$ws.Range("L7").Value = 1.025631846879809
$ws.Range("M7").Value = 1.038331505412127
$ws.Range("N7").Value = 1.01307202371974
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.021306071200111
$ws.Range("D8").Value = 1.036852283434038
$ws.Range("E8").Value = 1.02219614116701
$ws.Range("F8").Value = 1.034407803184389
$ws.Range("I8").Value = 1.035955456504621
$ws.Range("J8").Value = 1.026377679205355
$ws.Range("K8").Value = 1.039581421536885
$ws.Range("L8").Value = 1.024966797943518
$ws.Range("M8").Value = 1.03714377745012
$ws.Range("N8").Value = 1.012868549189484
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.019391422312717
$ws.Range("D9").Value = 1.035568518601576
$ws.Range("E9").Value = 1.020581857943989
$ws.Range("F9").Value = 1.031874791523706
$ws.Range("I9").Value = 1.03547803749739
$ws.Range("J9").Value = 1.025296708817325
$ws.Range("K9").Value = 1.038732504569814
$ws.Range("L9").Value = 1.023795565856398
$ws.Range("M9").Value = 1.035050885246921
$ws.Range("N9").Value = 1.012509216061472
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.018117213780572
$ws.Range("D10").Value = 1.034712052356164
$ws.Range("E10").Value = 1.019509690926737
$ws.Range("F10").Value = 1.030189202715895
$ws.Range("I10").Value = 1.035152862967282
$ws.Range("J10").Value = 1.024574960910227
$ws.Range("K10").Value = 1.038162258935866
$ws.Range("L10").Value = 1.023015260914925
$ws.Range("M10").Value = 1.033655691868395
$ws.Range("N10").Value = 1.012269154525859
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.017566011047986
$ws.Range("D11").Value = 1.034341067522891
$ws.Range("E11").Value = 1.019046404063131
$ws.Range("F11").Value = 1.029460063798771
$ws.Range("I11").Value = 1.035010432881533
$ws.Range("J11").Value = 1.024262187196772
$ws.Range("K11").Value = 1.037914331931663
$ws.Range("L11").Value = 1.022677516675463
$ws.Range("M11").Value = 1.033051582521949
$ws.Range("N11").Value = 1.012165089709134
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.017361351901742
$ws.Range("D12").Value = 1.034203249602791
$ws.Range("E12").Value = 1.018874465764092
$ws.Range("F12").Value = 1.029189339682342
$ws.Range("I12").Value = 1.034957284248851
$ws.Range("J12").Value = 1.02414597220076
$ws.Range("K12").Value = 1.037822090783439
$ws.Range("L12").Value = 1.022552084508754
$ws.Range("M12").Value = 1.032827192827492
$ws.Range("N12").Value = 1.012126418244136
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.017405248248203
$ws.Range("D13").Value = 1.034232812783831
$ws.Range("E13").Value = 1.018911340432674
$ws.Range("F13").Value = 1.029247405913273
$ws.Range("I13").Value = 1.034968695829687
$ws.Range("J13").Value = 1.024170902366821
$ws.Range("K13").Value = 1.037841883591501
$ws.Range("L13").Value = 1.022578989162541
$ws.Range("M13").Value = 1.032875325014131
$ws.Range("N13").Value = 1.012134714175383
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.017549092177151
$ws.Range("D14").Value = 1.034329675800042
$ws.Range("E14").Value = 1.019032188580509
$ws.Range("F14").Value = 1.029437683400667
$ws.Range("I14").Value = 1.035006044567292
$ws.Range("J14").Value = 1.024252581576386
$ws.Range("K14").Value = 1.037906710302708
$ws.Range("L14").Value = 1.022667147976972
$ws.Range("M14").Value = 1.033033034326542
$ws.Range("N14").Value = 1.012161893461738
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.017637730084993
$ws.Range("D15").Value = 1.034389354017845
$ws.Range("E15").Value = 1.019106666647655
$ws.Range("F15").Value = 1.029554934304857
$ws.Range("I15").Value = 1.035029024079023
$ws.Range("J15").Value = 1.02430290197068
$ws.Range("K15").Value = 1.037946632337337
$ws.Range("L15").Value = 1.022721468358804
$ws.Range("M15").Value = 1.033130204710113
$ws.Range("N15").Value = 1.012178637254092
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.018153806733129
$ws.Range("D16").Value = 1.034736670826597
$ws.Range("E16").Value = 1.019540458289074
$ws.Range("F16").Value = 1.030237608736539
$ws.Range("I16").Value = 1.035162281359103
$ws.Range("J16").Value = 1.024595713442175
$ws.Range("K16").Value = 1.03817869192779
$ws.Range("L16").Value = 1.02303767878346
$ws.Range("M16").Value = 1.033695785060845
$ws.Range("N16").Value = 1.012276058539171
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.018477672903959
$ws.Range("D17").Value = 1.03495450013449
$ws.Range("E17").Value = 1.019812824730042
$ws.Range("F17").Value = 1.030666028774188
$ws.Range("I17").Value = 1.035245434957867
$ws.Range("J17").Value = 1.024779319526798
$ws.Range("K17").Value = 1.038323988092785
$ws.Range("L17").Value = 1.023236065605702
$ws.Range("M17").Value = 1.03405056408828
$ws.Range("N17").Value = 1.012337137376175
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.018666630383346
$ws.Range("D18").Value = 1.035081543689423
$ws.Range("E18").Value = 1.019971784811738
$ws.Range("F18").Value = 1.030915989540012
$ws.Range("I18").Value = 1.035293779978391
$ws.Range("J18").Value = 1.02488638946039
$ws.Range("K18").Value = 1.038408639654358
$ws.Range("L18").Value = 1.023351794036205
$ws.Range("M18").Value = 1.034257502444763
$ws.Range("N18").Value = 1.012372752355485
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.018731068751981
$ws.Range("D19").Value = 1.035124860067647
$ws.Range("E19").Value = 1.020026001877094
$ws.Range("F19").Value = 1.031001231661141
$ws.Range("I19").Value = 1.035310237717366
$ws.Range("J19").Value = 1.024922893374641
$ws.Range("K19").Value = 1.038437487118643
$ws.Range("L19").Value = 1.023391256569801
$ws.Range("M19").Value = 1.034328063354895
$ws.Range("N19").Value = 1.012384894220084
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.018442919759573
$ws.Range("D20").Value = 1.034931130401086
$ws.Range("E20").Value = 1.019783592722222
$ws.Range("F20").Value = 1.030620056051541
$ws.Range("I20").Value = 1.035236529604501
$ws.Range("J20").Value = 1.024759622847263
$ws.Range("K20").Value = 1.038308409237583
$ws.Range("L20").Value = 1.02321477925002
$ws.Range("M20").Value = 1.034012499470438
$ws.Range("N20").Value = 1.012330585358436
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.017506731448999
$ws.Range("D21").Value = 1.034301152528482
$ws.Range("E21").Value = 1.018996597740678
$ws.Range("F21").Value = 1.029381648367277
$ws.Range("I21").Value = 1.034995053017939
$ws.Range("J21").Value = 1.024228530096205
$ws.Range("K21").Value = 1.037887624583266
$ws.Range("L21").Value = 1.022641186812573
$ws.Range("M21").Value = 1.032986592765569
$ws.Range("N21").Value = 1.012153890310127
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.016918588281909
$ws.Range("D22").Value = 1.033904959287096
$ws.Range("E22").Value = 1.018502634430601
$ws.Range("F22").Value = 1.028603652139561
$ws.Range("I22").Value = 1.034841817164099
$ws.Range("J22").Value = 1.023894397915872
$ws.Range("K22").Value = 1.037622193218315
$ws.Range("L22").Value = 1.022280669001472
$ws.Range("M22").Value = 1.032341584568701
$ws.Range("N22").Value = 1.012042696022412
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.017230328599854
$ws.Range("D23").Value = 1.034114997812366
$ws.Range("E23").Value = 1.018764412479156
$ws.Range("F23").Value = 1.029016021828541
$ws.Range("I23").Value = 1.034923183789838
$ws.Range("J23").Value = 1.024071547583578
$ws.Range("K23").Value = 1.037762985204756
$ws.Range("L23").Value = 1.022471774331532
$ws.Range("M23").Value = 1.032683513507544
$ws.Range("N23").Value = 1.012101651492372
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.018458623049594
$ws.Range("D24").Value = 1.034941690216807
$ws.Range("E24").Value = 1.019796801120045
$ws.Range("F24").Value = 1.030640828927825
$ws.Range("I24").Value = 1.035240554035743
$ws.Range("J24").Value = 1.024768523004368
$ws.Range("K24").Value = 1.038315448952241
$ws.Range("L24").Value = 1.023224397598609
$ws.Range("M24").Value = 1.034029699227333
$ws.Range("N24").Value = 1.012333545968093
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.019886018626009
$ws.Range("D25").Value = 1.035900520123147
$ws.Range("E25").Value = 1.020998486724425
$ws.Range("F25").Value = 1.03252909484129
$ws.Range("I25").Value = 1.035602680427703
$ws.Range("J25").Value = 1.025576363927272
$ws.Range("K25").Value = 1.038952733534086
$ws.Range("L25").Value = 1.024098271263794
$ws.Range("M25").Value = 1.035591939256273
$ws.Range("N25").Value = 1.012602203079232
